{"js": "// feat(django): add new placeholder responsible_person\n//\n// Insert a new paragraph \"Verantwortliche Person: {{ responsible_person }}\"\n// directly after the existing \"Publikationsdatum: {{ publication_date }}\"\n// paragraph, matching its paragraph/run formatting (Normal style, sz/szCs 24,\n// de-CH/zh-CN/hi-IN language run properties).\n\nconst body = context.document.body;\n\n// Locate the anchor paragraph robustly via search instead of a hard-coded\n// index, so the script keeps working even if earlier content shifts.\nconst results = body.search(\"Publikationsdatum:\", { matchCase: false, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor paragraph 'Publikationsdatum:' not found\");\n}\n\nconst anchorParagraph = results.items[0].paragraphs.getFirst();\nanchorParagraph.load(\"text\");\nawait context.sync();\n\n// Insert the new paragraph right after the anchor; it inherits the anchor's\n// paragraph/run formatting (Normal style, sz=24, szCs=24,\n// lang de-CH/zh-CN/hi-IN), which already matches the target.\nanchorParagraph.insertParagraph(\n  \"Verantwortliche Person: {{ responsible_person }}\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# feat(django): add new placeholder responsible_person\n#\n# Insert a new paragraph \"Verantwortliche Person: {{ responsible_person }}\"\n# directly after the existing \"Publikationsdatum: {{ publication_date }}\"\n# paragraph, matching its paragraph/run formatting (Normal style, sz/szCs 24,\n# de-CH/zh-CN/hi-IN language run properties).\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph robustly via Find instead of a hard-coded\n# index, so the script keeps working even if earlier content shifts.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Publikationsdatum:\")\n\nif (-not $found) {\n    throw \"Anchor paragraph 'Publikationsdatum:' not found\"\n}\n\n# Expand the found hit to its whole containing paragraph (wdParagraph = 4).\n[void]$rng.Expand(4)\n$anchorPara = $rng.Paragraphs(1)\n\n# Insert a new, empty paragraph right after the anchor; it inherits the\n# anchor's paragraph formatting (Normal / sz=24 / szCs=24), which already\n# matches the target.\n[void]$anchorPara.Range.InsertParagraphAfter()\n\n$insertedPara = $anchorPara.Next()\n$insertedPara.Range.Text = \"Verantwortliche Person: {{ responsible_person }}\"\n"}
